# Update "想去人数" (interest count) values in both the "展览" and "全部类型"
# worksheets. Both sheets contain the same rows of data, so the same set of
# cell updates is applied to each sheet.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F2"  = 149
    "F4"  = 12113
    "F5"  = 1252
    "F6"  = 131
    "F7"  = 26
    "F8"  = 88
    "F10" = 187
    "F11" = 438
    "F16" = 352
    "F17" = 2297
    "F18" = 87
    "F19" = 928
    "F20" = 120
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
